$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 281, pushing all
# existing rows from 281-298 down to 282-299 (the rest of the sheet is
# unchanged, just shifted down by one row).
$ws.Rows("281").Insert()

# Populate the newly inserted row 281 with the new record's data.
$ws.Range("A281").Value = 4
$ws.Range("B281").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C281").Value = 'Los Lagos'
$ws.Range("D281").Value = 44826
$ws.Range("D281").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = 100112032
$ws.Range("G281").Value = 'Zapallo italiano'
$ws.Range("H281").Value = 'Sin especificar'
$ws.Range("I281").Value = 'Primera'
$ws.Range("J281").Value = 120
$ws.Range("K281").Value = 22000
$ws.Range("L281").Value = 22000
$ws.Range("M281").Value = 22000
$ws.Range("N281").Value = '$/caja 50 unidades'
$ws.Range("O281").Value = 'Región de Arica y Parinacota'
$ws.Range("P281").Value = 440
$ws.Range("Q281").Value = 50
$ws.Range("R281").Value = 'Hortaliza'
